$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Vega Modelo de Temuco" (Arándano).
# It belongs chronologically at row 26, pushing the existing data rows
# (old rows 26..118) down by one (new rows 27..119).
$ws.Rows("26:26").Insert()

# Populate the newly inserted row 26 with the new record's data.
$ws.Range("A26").Value = 10
$ws.Range("B26").Value = "Vega Modelo de Temuco"
$ws.Range("C26").Value = "La Araucanía"
$ws.Range("D26").Value = 44910
$ws.Range("E26").Value = 9
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100101
$ws.Range("H26").Value = "Berries"
$ws.Range("I26").Value = 100101001
$ws.Range("J26").Value = "Arándano (blue)"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 400
$ws.Range("N26").Value = 2000
$ws.Range("O26").Value = 2000
$ws.Range("P26").Value = 2000
$ws.Range("Q26").Value = "$/kilo"
$ws.Range("R26").Value = "Región del Maule"
$ws.Range("S26").Value = 2000
$ws.Range("T26").Value = 1
